$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.022.36"
$ws.Range("E2").Value = "  +1.49%  "
$ws.Range("D3").Value = "2.418.66"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'555.87"
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("D6").Value = "'143.00"
$ws.Range("E6").Value = "  +3.38%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.532"
$ws.Range("E8").Value = "  +1.71%  "
$ws.Range("D9").Value = "2.416.16"
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("D11").Value = "'0.156"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "'5.39"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "'0.353"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").Value = "'26.20"
$ws.Range("E14").Value = "  +5.23%  "
$ws.Range("D15").Value = "'0.0000176"
$ws.Range("E15").Value = "  +8.21%  "
$ws.Range("D16").Value = "2.855.75"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("D17").Value = "62.015.58"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "2.415.39"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").Value = "'11.09"
$ws.Range("E19").Value = "  +3.17%  "
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").Value = "'323.92"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'65.31"
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("B25").Value = "SuiNetwork"
$ws.Range("C25").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D25").Value = "'1.73"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").Value = "'8.96"
$ws.Range("E26").Value = "  +5.01%  "
$ws.Range("D27").Value = "'574.10"
$ws.Range("E27").Value = "  +14.11%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.529.54"
$ws.Range("E29").Value = "  +2.31%  "
$ws.Range("D30").Value = "0.0₃0938"
$ws.Range("E30").Value = "  +5.97%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.47"
$ws.Range("E31").Value = "  +5.64%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'8.28"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").Value = "'1.87"
$ws.Range("E34").Value = "  +2.89%  "
$ws.Range("E35").Value = "  +3.30%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").Value = "'5.67"
$ws.Range("E37").Value = "  +6.14%  "
$ws.Range("D38").Value = "'4.81"
$ws.Range("E38").Value = "  +2.30%  "
$ws.Range("D39").Value = "'0.383"
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "'18.75"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.85"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("D42").Value = "'149.17"
$ws.Range("E42").Value = "  +3.07%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.31"
$ws.Range("E44").Value = "  +9.86%  "
$ws.Range("D45").Value = "'149.99"
$ws.Range("E45").Value = "  +4.38%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'3.65"
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0546"
$ws.Range("E47").Value = "  +4.51%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'20.31"
$ws.Range("E48").Value = "  +5.25%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.591"
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.0918"
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0228"
$ws.Range("E51").Value = "  +2.67%  "
